$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5 (shifts existing rows 5-24 down to 6-25,
# carrying formatting along) to make room for a new weekly price record.
$ws.Rows("5:5").Insert()

# Populate the newly inserted row 5 with the new record.
$ws.Range("A5").Value = 6
$ws.Range("B5").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C5").Value = "Metropolitana"
$ws.Range("D5").Value = 44600
$ws.Range("E5").Value = 13
$ws.Range("F5").Value = 300000001
$ws.Range("G5").Value = "Rabanito"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 1300
$ws.Range("K5").Value = 3500
$ws.Range("L5").Value = 4000
$ws.Range("M5").Value = 3808
$ws.Range("N5").Value = '$/cien unidades (volumen en unidades)'
$ws.Range("O5").Value = "Región Metropolitana"
$ws.Range("P5").Value = 38
$ws.Range("Q5").Value = 100
$ws.Range("R5").Value = "Hortaliza"
